$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H: tidal elevation converted to feet above a fixed datum (G + 16 inches).
# H9 is entered as its own formula (mirrors how G9 stands alone before the
# shared formula block starts at G10), and H10:H32 is entered as one range so
# Excel records it as a shared formula, mirroring column G's G10:G32 pattern.
$ws.Range("H9").Formula = "=G9+16/12"
$ws.Range("H10:H32").Formula = "=G10+16/12"

# Apply the custom 6-decimal number format to the whole new column.
$ws.Range("H9:H32").NumberFormat = "0.000000"

# Match the recorded selection state from the edit.
$ws.Range("H9:H32").Select()
